$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: grow the two 11-column stat blocks (B:L and M:W) to 13 columns
# each (B:N and O:AA), adding room for the new "sum" and "msg_per_mus"
# columns, while keeping the row-1 group labels ("msg_count_twitter" /
# "msg_count_facebook") anchored at B1 / (shifted) O1 respectively.
#
# Inserting a column at a position strictly inside an existing merged cell
# grows the merge instead of shifting its anchor, so we insert at C (inside
# B1:L1) twice, and then inside the (now shifted) second block twice as well.

$ws.Range("C1").EntireColumn.Insert()
$ws.Range("C1").EntireColumn.Insert()
$ws.Range("P1").EntireColumn.Insert()
$ws.Range("P1").EntireColumn.Insert()

# --- Step 2: write the final header row (row 2) for both blocks.
$headers = @("sum","mean","std","min","q25","median","q75","max","count","msg_per_mus","active_mus_n","active_mus_pc","active_mus_pc_z")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(2, 2 + $i).Value = $headers[$i]
    $ws.Cells.Item(2, 15 + $i).Value = $headers[$i]
}

# --- Step 3: write the data rows (row 4 = Accredited, row 5 = Unaccredited)
# for both blocks.
$row4_block1 = @(3321036, 2011.5, 5267.7, 0, 1, 471, 1765, 75337, 1651, 2011.5, 1240, 75.09999999999999, 0.7)
$row4_block2 = @(766527, 464.3, 1048.9, 0, 0, 0, 615, 13530, 1651, 464.3, 731, 44.3, 0.7)
$row5_block1 = @(2292591, 1358.2, 4794.8, 0, 0, 19.5, 737.5, 74189, 1688, 1358.2, 959, 56.8, -0.7)
$row5_block2 = @(701065, 415.3, 1896.8, 0, 0, 0, 241.5, 37660, 1688, 415.3, 518, 30.7, -0.7)

for ($i = 0; $i -lt $row4_block1.Length; $i++) {
    $ws.Cells.Item(4, 2 + $i).Value = $row4_block1[$i]
    $ws.Cells.Item(4, 15 + $i).Value = $row4_block2[$i]
    $ws.Cells.Item(5, 2 + $i).Value = $row5_block1[$i]
    $ws.Cells.Item(5, 15 + $i).Value = $row5_block2[$i]
}
